# Fruta / hortaliza, semanal
# Inserts two new weekly price records (Provincia de Los Andes, 2022-02-18)
# above the existing "Femacal de La Calera - Tuna" history, pushing the
# remaining rows down by two positions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 97-98; everything currently at row 97 onward
# shifts down to make room.
$ws.Rows("97:98").Insert()

# New row 97: Primera quality record for Provincia de Los Andes.
$ws.Range("A97").Value = 3
$ws.Range("B97").Value = "Femacal de La Calera"
$ws.Range("C97").Value = "Coquimbo"
$ws.Range("D97").Value = 44610
$ws.Range("E97").Value = 5
$ws.Range("F97").Value = "Fruta"
$ws.Range("G97").Value = 100107
$ws.Range("H97").Value = "Otros"
$ws.Range("I97").Value = 100107011
$ws.Range("J97").Value = "Tuna"
$ws.Range("K97").Value = "Sin especificar"
$ws.Range("L97").Value = "Primera"
$ws.Range("M97").Value = 50
$ws.Range("N97").Value = 16000
$ws.Range("O97").Value = 16000
$ws.Range("P97").Value = 16000
$ws.Range("Q97").Value = "$/caja 16 kilos"
$ws.Range("R97").Value = "Provincia de Los Andes"
$ws.Range("S97").Value = 1000
$ws.Range("T97").Value = 16

# New row 98: Segunda quality record for Provincia de Los Andes.
$ws.Range("A98").Value = 3
$ws.Range("B98").Value = "Femacal de La Calera"
$ws.Range("C98").Value = "Coquimbo"
$ws.Range("D98").Value = 44610
$ws.Range("E98").Value = 5
$ws.Range("F98").Value = "Fruta"
$ws.Range("G98").Value = 100107
$ws.Range("H98").Value = "Otros"
$ws.Range("I98").Value = 100107011
$ws.Range("J98").Value = "Tuna"
$ws.Range("K98").Value = "Sin especificar"
$ws.Range("L98").Value = "Segunda"
$ws.Range("M98").Value = 50
$ws.Range("N98").Value = 14000
$ws.Range("O98").Value = 14000
$ws.Range("P98").Value = 14000
$ws.Range("Q98").Value = "$/caja 16 kilos"
$ws.Range("R98").Value = "Provincia de Los Andes"
$ws.Range("S98").Value = 875
$ws.Range("T98").Value = 16
